$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "survey area"
$ws.Range("B2").Value = 319812.5
$ws.Range("C2").Value = "6205 (815 - 13819)"
$ws.Range("D2").Value = "0.0024 (3e-04 - 0.0054)"
$ws.Range("E2").Value = "0.0194 (0.0025 - 0.0432)"

$ws.Range("A3").Value = "hotspots"
$ws.Range("B3").Value = 20381.25
$ws.Range("C3").Value = "3492 (801 - 6640)"
$ws.Range("D3").Value = "0.0213 (0.0049 - 0.0404)"
$ws.Range("E3").Value = "0.1714 (0.0393 - 0.3258)"
